$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 725.8182
$ws.Range("I53").Value = 290.4
$ws.Range("J53").Value = 1088.6666
$ws.Range("K53").Value = 290.4
$ws.Range("L53").Value = 1088.6666
$ws.Range("M53").Value = 346.6
$ws.Range("N53").Value = -2362.6666
$ws.Range("H64").Value = 6035.5713
$ws.Range("I64").Value = 4291.5
$ws.Range("K64").Value = 4291.5
$ws.Range("M64").Value = -4043.5
$ws.Range("H67").Value = 6035.5713
$ws.Range("I67").Value = 4291.5
$ws.Range("K67").Value = 4291.5
$ws.Range("M67").Value = -3433.5
$ws.Range("H101").Value = 336.81818
$ws.Range("I101").Value = 336.81818
$ws.Range("K101").Value = 1010.45454
$ws.Range("M101").Value = 611.54546
$ws.Range("H107").Value = 1043.5
$ws.Range("I107").Value = 1038.55
$ws.Range("J107").Value = 1093
$ws.Range("K107").Value = 1038.55
$ws.Range("L107").Value = 1093
$ws.Range("M107").Value = 881.45
$ws.Range("N107").Value = -4933
$ws.Range("H125").Value = 8120.8
$ws.Range("J125").Value = 7499.375
$ws.Range("L125").Value = 67494.375
$ws.Range("N125").Value = -72414.375
$ws.Range("H135").Value = 1043.56
$ws.Range("J135").Value = 2888
$ws.Range("L135").Value = 25992
$ws.Range("N135").Value = -31062
$ws.Range("H137").Value = 3734.9285
$ws.Range("I137").Value = 4154.4443
$ws.Range("K137").Value = 12463.3329
$ws.Range("M137").Value = -9913.332900000001
$ws.Range("H138").Value = 9622753
$ws.Range("J138").Value = 10211642
$ws.Range("L138").Value = 30634926
$ws.Range("N138").Value = -30645206

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22981.611
$ws.Range("I32").Value = 21266.158
$ws.Range("K32").Value = 21266.158
$ws.Range("M32").Value = -20979.158
$ws.Range("H61").Value = 262656.53
$ws.Range("J61").Value = 733648.6
$ws.Range("L61").Value = 733648.6
$ws.Range("N61").Value = -734072.6
$ws.Range("H74").Value = 28352.666
$ws.Range("I74").Value = 14136.889
$ws.Range("K74").Value = 14136.889
$ws.Range("M74").Value = -13262.889
$ws.Range("H77").Value = 28352.666
$ws.Range("I77").Value = 14136.889
$ws.Range("K77").Value = 70684.44499999999
$ws.Range("M77").Value = -66316.44499999999
$ws.Range("H110").Value = 12411.95
$ws.Range("I110").Value = 15855
$ws.Range("K110").Value = 15855
$ws.Range("M110").Value = -13810
$ws.Range("H122").Value = 1630.2142
$ws.Range("I122").Value = 1652.3334
$ws.Range("K122").Value = 4957.0002
$ws.Range("M122").Value = -2507.0002
$ws.Range("H132").Value = 3750
$ws.Range("I132").Value = 4250
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 12750
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -10220
$ws.Range("N132").Value = -14810
$ws.Range("H136").Value = 262656.53
$ws.Range("J136").Value = 733648.6
$ws.Range("L136").Value = 2200945.8
$ws.Range("N136").Value = -2206045.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3119.818
$ws.Range("I20").Value = 3304.7
$ws.Range("K20").Value = 3304.7
$ws.Range("M20").Value = -3057.7
$ws.Range("H107").Value = 1037.762
$ws.Range("I107").Value = 1055.4736
$ws.Range("K107").Value = 1055.4736
$ws.Range("M107").Value = 864.5264
$ws.Range("H134").Value = 2030.1428
$ws.Range("I134").Value = 2030.1428
$ws.Range("K134").Value = 6090.428400000001
$ws.Range("M134").Value = -3555.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3791.2354
$ws.Range("I31").Value = 3590.6875
$ws.Range("K31").Value = 3590.6875
$ws.Range("M31").Value = -3295.6875
$ws.Range("H34").Value = 3791.2354
$ws.Range("I34").Value = 3590.6875
$ws.Range("K34").Value = 3590.6875
$ws.Range("M34").Value = -3388.6875
$ws.Range("H62").Value = 9451.25
$ws.Range("I62").Value = 7601.6665
$ws.Range("K62").Value = 7601.6665
$ws.Range("M62").Value = -6977.6665
$ws.Range("H65").Value = 9451.25
$ws.Range("I65").Value = 7601.6665
$ws.Range("K65").Value = 38008.3325
$ws.Range("M65").Value = -34888.3325
$ws.Range("H87").Value = 75000
$ws.Range("J87").Value = 75000
$ws.Range("L87").Value = 75000
$ws.Range("N87").Value = -77372
$ws.Range("H90").Value = 75000
$ws.Range("J90").Value = 75000
$ws.Range("L90").Value = 225000
$ws.Range("N90").Value = -236856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1684190.5
$ws.Range("J34").Value = 500.5
$ws.Range("L34").Value = 1501.5
$ws.Range("N34").Value = -1669.5
$ws.Range("H39").Value = 1333.6666
$ws.Range("J39").Value = 1333.6666
$ws.Range("L39").Value = 4000.9998
$ws.Range("N39").Value = -4588.9998
$ws.Range("H55").Value = 1456.6666
$ws.Range("J55").Value = 1600.2
$ws.Range("L55").Value = 4800.6
$ws.Range("N55").Value = -5154.6
$ws.Range("H131").Value = 3712.5833
$ws.Range("J131").Value = 4823.5884
$ws.Range("L131").Value = 14470.7652
$ws.Range("N131").Value = -24550.7652

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 13930.238
$ws.Range("J126").Value = 4282.6665
$ws.Range("L126").Value = 12847.9995
$ws.Range("N126").Value = -17787.9995
$ws.Range("H132").Value = 6108.1665
$ws.Range("I132").Value = 3900
$ws.Range("J132").Value = 7212.25
$ws.Range("K132").Value = 11700
$ws.Range("L132").Value = 21636.75
$ws.Range("M132").Value = -9170
$ws.Range("N132").Value = -26696.75
$ws.Range("H133").Value = 95195
$ws.Range("J133").Value = 95195
$ws.Range("L133").Value = 95195
$ws.Range("N133").Value = -105315
$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140
$ws.Range("H138").Value = 84674
$ws.Range("J138").Value = 84674
$ws.Range("L138").Value = 84674
$ws.Range("N138").Value = -94954
$ws.Range("H139").Value = 103098.664
$ws.Range("J139").Value = 109500
$ws.Range("L139").Value = 109500
$ws.Range("N139").Value = -119780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4443.1665
$ws.Range("I7").Value = 4364.9165
$ws.Range("K7").Value = 4364.9165
$ws.Range("M7").Value = -4252.9165
$ws.Range("H20").Value = 22400
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20452
$ws.Range("H61").Value = 63505.41
$ws.Range("I61").Value = 87174.336
$ws.Range("K61").Value = 87174.336
$ws.Range("M61").Value = -86972.336
$ws.Range("H113").Value = 63505.41
$ws.Range("I113").Value = 87174.336
$ws.Range("K113").Value = 87174.336
$ws.Range("M113").Value = -85004.336
$ws.Range("H122").Value = 9399.6
$ws.Range("I122").Value = 9000
$ws.Range("J122").Value = 9666
$ws.Range("K122").Value = 27000
$ws.Range("L122").Value = 28998
$ws.Range("M122").Value = -24550
$ws.Range("N122").Value = -33898
$ws.Range("H126").Value = 4443.1665
$ws.Range("I126").Value = 4364.9165
$ws.Range("K126").Value = 13094.7495
$ws.Range("M126").Value = -10624.7495
$ws.Range("H136").Value = 4704.25
$ws.Range("I136").Value = 4581.514
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 13744.542
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -11194.542
$ws.Range("N136").Value = -32100
$ws.Range("H138").Value = 90404.664
$ws.Range("J138").Value = 90404.664
$ws.Range("L138").Value = 90404.664
$ws.Range("N138").Value = -100684.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15540.714
$ws.Range("J62").Value = 21928.334
$ws.Range("L62").Value = 21928.334
$ws.Range("N62").Value = -23176.334
$ws.Range("H65").Value = 15540.714
$ws.Range("J65").Value = 21928.334
$ws.Range("L65").Value = 109641.67
$ws.Range("N65").Value = -115881.67
$ws.Range("H132").Value = 4268.923
$ws.Range("I132").Value = 4117.353
$ws.Range("K132").Value = 12352.059
$ws.Range("M132").Value = -9822.059000000001
$ws.Range("H140").Value = 147998.56
$ws.Range("J140").Value = 72748.375
$ws.Range("L140").Value = 72748.375
$ws.Range("N140").Value = -83108.375
